# Remove the staff record for "Dr. ISWARYA M" (row 14), which shifts the
# remaining rows (Vijayakumar, Antony Gnana Aravind, Jerlin Sheeba, Muthuvel)
# up by one row, and updates the sheet's dimension/used-range accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Delete()

# Restore the selected cell as recorded in the saved workbook.
$ws.Range("A21").Select()
